# Weekly fruit/vegetable price update: a new price-survey row for
# "Berenjena" (Mercado Mayorista Lo Valledor de Santiago) is inserted
# right before the existing row 199, pushing the old rows 199-236 down
# to 200-237 and extending the used range from A1:R236 to A1:R237.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 199 (shifts rows 199..236 down to 200..237).
$ws.Rows.Item(199).Insert()

# Populate the newly inserted row 199 with the new survey record.
$ws.Cells.Item(199, 1).Value  = 6
$ws.Cells.Item(199, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(199, 3).Value  = "Metropolitana"
$ws.Cells.Item(199, 4).Value  = 44798
$ws.Cells.Item(199, 5).Value  = 13
$ws.Cells.Item(199, 6).Value  = 100112001
$ws.Cells.Item(199, 7).Value  = "Berenjena"
$ws.Cells.Item(199, 8).Value  = "Sin especificar"
$ws.Cells.Item(199, 9).Value  = "Primera"
$ws.Cells.Item(199, 10).Value = 580
$ws.Cells.Item(199, 11).Value = 10000
$ws.Cells.Item(199, 12).Value = 11000
$ws.Cells.Item(199, 13).Value = 10448
$ws.Cells.Item(199, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(199, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(199, 16).Value = 261
$ws.Cells.Item(199, 17).Value = 40
$ws.Cells.Item(199, 18).Value = "Hortaliza"

# Give the new date cell the same date style as the rest of column D.
$ws.Cells.Item(199, 4).NumberFormat = $ws.Cells.Item(200, 4).NumberFormat
